$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.51%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'32.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'2.62%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.144"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.95%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07867"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.75%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.279"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-6.02%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.817"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.27%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.817"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.50%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9321"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.11%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1777"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.04%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07703"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'5.28%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08854"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.19%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03081"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'1.54%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.1006"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.64%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001513"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.71%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005983"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.14%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.465"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.03%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.251"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.22%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3271"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.43%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1339"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.17%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.284"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-6.58%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.1810"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'11.97%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04620"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.25%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001254"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.79%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004521"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'2.04%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001253"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'4.44%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-1.29%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01788"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'0.12%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04732"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'6.00%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007251"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'4.60%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1383"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'2.80%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002126"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-3.77%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01009"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'2.96%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006346"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-3.40%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.21%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.003206"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-38.78%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.7333"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-10.64%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002104"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.21%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002004"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.21%"
$ws.Range("E50").Style = "Normal"
